# Refresh crypto price/volume data to match the latest coinranking.com snapshot.
# Also reflects two ranking swaps that occurred between pulls:
#   - Uniswap and PEPE traded places (rows 23/24)
#   - OKB, ImmutableX and VeChain reshuffled (rows 49/50/51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "92.203.75"
$ws.Range("E2").Value = "  +0.76%  "

$ws.Range("D3").Value = "3.100.40"
$ws.Range("E3").Value = "  -0.92%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'233.53"
$ws.Range("E5").Value = "  -2.61%  "

$ws.Range("D6").Value = "'613.50"
$ws.Range("E6").Value = "  -0.60%  "

$ws.Range("D7").Value = "'1.08"
$ws.Range("E7").Value = "  -2.58%  "

$ws.Range("D8").Value = "'0.386"
$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("D10").Value = "3.095.63"
$ws.Range("E10").Value = "  -0.96%  "

$ws.Range("E11").Value = "  +5.68%  "

$ws.Range("E12").Value = "  -2.92%  "

$ws.Range("E13").Value = "  -4.42%  "

$ws.Range("D14").Value = "92.083.93"
$ws.Range("E14").Value = "  +0.84%  "

$ws.Range("D15").Value = "'33.80"
$ws.Range("E15").Value = "  -2.90%  "

$ws.Range("E16").Value = "  -2.95%  "

$ws.Range("D17").Value = "3.681.74"
$ws.Range("E17").Value = "  -1.36%  "

$ws.Range("D18").Value = "3.093.57"
$ws.Range("E18").Value = "  -1.23%  "

$ws.Range("D19").Value = "'3.81"
$ws.Range("E19").Value = "  +1.99%  "

$ws.Range("D20").Value = "'14.38"
$ws.Range("E20").Value = "  -3.24%  "

$ws.Range("E21").Value = "  -1.95%  "

$ws.Range("D22").Value = "'437.08"
$ws.Range("E22").Value = "  -3.92%  "

$ws.Range("B23").Value = "PEPE"
$ws.Range("C23").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D23").Value = "'0.0000199"
$ws.Range("E23").Value = "  -1.20%  "

$ws.Range("B24").Value = "Uniswap"
$ws.Range("C24").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D24").Value = "'9.11"
$ws.Range("E24").Value = "  -0.42%  "

$ws.Range("E25").Value = "  -1.84%  "

$ws.Range("D26").Value = "'85.22"
$ws.Range("E26").Value = "  -3.54%  "

$ws.Range("D27").Value = "'11.37"
$ws.Range("E27").Value = "  -3.15%  "

$ws.Range("D28").Value = "3.259.35"

$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("E30").Value = "  +7.41%  "

$ws.Range("E31").Value = "  +3.07%  "

$ws.Range("E32").Value = "  -18.32%  "

$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +31.63%  "

$ws.Range("D34").Value = "'9.11"
$ws.Range("E34").Value = "  -2.42%  "

$ws.Range("D35").Value = "'8.05"
$ws.Range("E35").Value = "  +8.55%  "

$ws.Range("D36").Value = "'0.156"
$ws.Range("E36").Value = "  -10.74%  "

$ws.Range("E37").Value = "  -2.46%  "

$ws.Range("E38").Value = "  -0.61%  "

$ws.Range("D39").Value = "'1.88"
$ws.Range("E39").Value = "  -5.17%  "

$ws.Range("D40").Value = "'23.85"
$ws.Range("E40").Value = "  +7.70%  "

$ws.Range("E41").Value = "  -0.34%  "

$ws.Range("E42").Value = "  -3.31%  "

$ws.Range("D43").Value = "'465.42"
$ws.Range("E43").Value = "  -5.13%  "

$ws.Range("D44").Value = "'3.27"
$ws.Range("E44").Value = "  -3.15%  "

$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("D46").Value = "'158.86"
$ws.Range("E46").Value = "  +1.70%  "

$ws.Range("E47").Value = "  -3.24%  "

$ws.Range("E48").Value = "  -4.66%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0326"
$ws.Range("E49").Value = "  +1.30%  "

$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").Value = "'43.74"
$ws.Range("E50").Value = "  -0.74%  "

$ws.Range("B51").Value = "ImmutableX"
$ws.Range("C51").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D51").Value = "'1.31"
$ws.Range("E51").Value = "  -2.70%  "
